$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price/volume data scraped on Mon Nov 11 08:56:06 UTC 2024

$ws.Range('D2').Value = '81.278.30'
$ws.Range('E2').Value = '  +2.78%  '
$ws.Range('D3').Value = '3.146.78'
$ws.Range('E3').Value = '  -1.12%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '208.97'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '621.72'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.283'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +23.13%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.582'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').Value = '3.145.16'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.584'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000251'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +11.19%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.29'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.87%  '
$ws.Range('D15').Value = '3.721.30'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.39'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.36%  '
$ws.Range('D17').Value = '81.127.58'
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').Value = '3.143.65'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.17'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +5.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.90'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '432.35'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.31%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '8.94'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('E23').Value = '  +2.38%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.19'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.24'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.40%  '
$ws.Range('D26').Value = '3.307.26'
$ws.Range('E26').Value = '  -1.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '76.19'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.81'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -3.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.09%  '
$ws.Range('E30').Value = '  +5.54%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '582.94'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +11.88%  '
$ws.Range('E32').Value = '  +0.77%  '
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.49'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.58%  '
$ws.Range('E35').Value = '  +11.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.139'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.99'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '22.71'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.10%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.999'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.01'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +10.49%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.04'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +14.02%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.05'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +21.42%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '20.72'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.63%  '
$ws.Range('E45').Value = '  -3.31%  '
$ws.Range('E46').Value = '  +0.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '186.31'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.84%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '45.29'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.40%  '
$ws.Range('E49').Value = '  +0.29%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.766'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.69%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.91'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.60%  '
